$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2654603333333334
$ws.Range("H2").Value = 0.796381
$ws.Range("I2").Value = 0.04079010536687974
$ws.Range("J2").Value = 0.04079010536687975
$ws.Range("M2").Value = 1.923239
$ws.Range("N2").Value = 5.769717
$ws.Range("O2").Value = 0.2340262838603868
$ws.Range("P2").Value = 0.2340262838603868
$ws.Range("Q2").Value = 0.5105436660196667
$ws.Range("R2").Value = 4.594892994177
$ws.Range("S2").Value = 0.009545956777284486
$ws.Range("T2").Value = 0.009545956777284486
$ws.Range("G3").Value = 0.2654603333333334
$ws.Range("H3").Value = 0.796381
$ws.Range("I3").Value = 0.04079010536687974
$ws.Range("J3").Value = 0.04079010536687975
$ws.Range("O3").Value = 0.4335574295612247
$ws.Range("P3").Value = 0.4335574295612246
$ws.Range("Q3").Value = 0.9458339288517779
$ws.Range("R3").Value = 8.512505359666
$ws.Range("S3").Value = 0.0176848532343959
$ws.Range("T3").Value = 0.0176848532343959
$ws.Range("G4").Value = 0.2654603333333334
$ws.Range("H4").Value = 0.796381
$ws.Range("I4").Value = 0.04079010536687974
$ws.Range("J4").Value = 0.04079010536687975
$ws.Range("O4").Value = 0.3324162865783886
$ws.Range("P4").Value = 0.3324162865783886
$ws.Range("Q4").Value = 0.7251879010975556
$ws.Range("R4").Value = 6.526691109878
$ws.Range("S4").Value = 0.01355929535519936
$ws.Range("T4").Value = 0.01355929535519936
$ws.Range("I5").Value = 0.8420553458721338
$ws.Range("J5").Value = 0.8420553458721339
$ws.Range("M5").Value = 1.923239
$ws.Range("N5").Value = 5.769717
$ws.Range("O5").Value = 0.2340262838603868
$ws.Range("P5").Value = 0.2340262838603868
$ws.Range("Q5").Value = 10.53946831973833
$ws.Range("R5").Value = 94.855214877645
$ws.Range("S5").Value = 0.1970630833992282
$ws.Range("T5").Value = 0.1970630833992282
$ws.Range("I6").Value = 0.8420553458721338
$ws.Range("J6").Value = 0.8420553458721339
$ws.Range("O6").Value = 0.4335574295612247
$ws.Range("P6").Value = 0.4335574295612246
$ws.Range("S6").Value = 0.3650793513046103
$ws.Range("T6").Value = 0.3650793513046103
$ws.Range("I7").Value = 0.8420553458721338
$ws.Range("J7").Value = 0.8420553458721339
$ws.Range("O7").Value = 0.3324162865783886
$ws.Range("P7").Value = 0.3324162865783886
$ws.Range("S7").Value = 0.2799129111682954
$ws.Range("T7").Value = 0.2799129111682954
$ws.Range("G8").Value = 0.7624369999999999
$ws.Range("I8").Value = 0.1171545487609863
$ws.Range("J8").Value = 0.1171545487609864
$ws.Range("M8").Value = 1.923239
$ws.Range("N8").Value = 5.769717
$ws.Range("O8").Value = 0.2340262838603868
$ws.Range("P8").Value = 0.2340262838603868
$ws.Range("Q8").Value = 1.466348573443
$ws.Range("R8").Value = 13.197137160987
$ws.Range("S8").Value = 0.02741724368387412
$ws.Range("T8").Value = 0.02741724368387412
$ws.Range("G9").Value = 0.7624369999999999
$ws.Range("I9").Value = 0.1171545487609863
$ws.Range("J9").Value = 0.1171545487609864
$ws.Range("O9").Value = 0.4335574295612247
$ws.Range("P9").Value = 0.4335574295612246
$ws.Range("Q9").Value = 2.716559472960666
$ws.Range("S9").Value = 0.0507932250222184
$ws.Range("T9").Value = 0.0507932250222184
$ws.Range("G10").Value = 0.7624369999999999
$ws.Range("I10").Value = 0.1171545487609863
$ws.Range("J10").Value = 0.1171545487609864
$ws.Range("O10").Value = 0.3324162865783886
$ws.Range("P10").Value = 0.3324162865783886
$ws.Range("Q10").Value = 2.082835054135333
$ws.Range("S10").Value = 0.03894408005489384
$ws.Range("T10").Value = 0.03894408005489384
